$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# ThemeColorScheme item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
# 11 hlink, 12 folHlink
$cs.Item(1).RGB = 0
$cs.Item(2).RGB = 16777215
$cs.Item(3).RGB = 6968388
$cs.Item(4).RGB = 15132391
$cs.Item(5).RGB = 13998939
$cs.Item(6).RGB = 3243501
$cs.Item(7).RGB = 10855845
$cs.Item(8).RGB = 49407
$cs.Item(9).RGB = 12874308
$cs.Item(10).RGB = 4697456
$cs.Item(11).RGB = 12673797
$cs.Item(12).RGB = 7491477
